$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Lookup" / "value" / "unit" header cells entirely (row 1)
$ws.Range("A1").Clear()
$ws.Range("C1").Clear()
$ws.Range("D1").Clear()

# Relabel the remaining row 1 headers
$ws.Range("B1").Value = "label"
$ws.Range("E1").Value = "description"

# Relabel row 2 key
$ws.Range("A2").Value = "report_type"

# Move the active selection from A16 to A2
$ws.Range("A2").Select()
